# Updates cryptos list price/volume data (and the Aave/BabyDogeCoin row swap)
# to match the refreshed scrape, per commit "Updated cryptos list on
# Mon Apr 17 17:09:20 UTC 2023 with GitHub Actions".

function Set-TextValue {
    # Forces a numeric-looking literal (e.g. "1.008") to be written as TEXT,
    # matching the source data's inline-string cells, instead of letting
    # Excel's smart-parse coerce it into a real number.
    param($Ws, $CellRef, $Val)
    $range = $Ws.Range($CellRef)
    $range.NumberFormat = "@"
    $range.Value = $Val
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.674.94'
$ws.Range('E2').Value = '  -3.05%  '
$ws.Range('D3').Value = '2.095.43'
$ws.Range('E3').Value = '  -1.27%  '
Set-TextValue $ws 'D4' '1.008'
$ws.Range('E4').Value = '  -0.39%  '
Set-TextValue $ws 'D5' '342.72'
$ws.Range('E5').Value = '  -2.07%  '
Set-TextValue $ws 'D7' '0.5122'
$ws.Range('E7').Value = '  -2.74%  '
Set-TextValue $ws 'D8' '0.4401'
$ws.Range('E8').Value = '  -2.43%  '
Set-TextValue $ws 'D9' '53.35'
$ws.Range('E9').Value = '  -1.66%  '
Set-TextValue $ws 'D10' '0.09140'
$ws.Range('E10').Value = '  +0.53%  '
Set-TextValue $ws 'D11' '1.170'
$ws.Range('E11').Value = '  -0.95%  '
Set-TextValue $ws 'D12' '24.77'
$ws.Range('E12').Value = '  +0.86%  '
$ws.Range('D13').Value = '2.106.06'
$ws.Range('E13').Value = '  -0.08%  '
Set-TextValue $ws 'D14' '6.735'
$ws.Range('E14').Value = '  -1.76%  '
Set-TextValue $ws 'D15' '8.195'
$ws.Range('E15').Value = '  +1.18%  '
Set-TextValue $ws 'D16' '99.56'
$ws.Range('E16').Value = '  -2.87%  '
Set-TextValue $ws 'D17' '0.00001149'
$ws.Range('E17').Value = '  -2.77%  '
$ws.Range('E19').Value = '  +8.13%  '
Set-TextValue $ws 'D20' '0.06644'
$ws.Range('E20').Value = '  -1.15%  '
$ws.Range('E21').Value = '  -0.30%  '
Set-TextValue $ws 'D22' '6.169'
$ws.Range('D23').Value = '29.735.39'
$ws.Range('E23').Value = '  -3.06%  '
$ws.Range('E24').Value = '  -2.06%  '
$ws.Range('E25').Value = '  -3.13%  '
$ws.Range('D26').Value = '2.348.35'
$ws.Range('E26').Value = '  -0.32%  '
Set-TextValue $ws 'D27' '21.82'
$ws.Range('E27').Value = '  -3.06%  '
Set-TextValue $ws 'D28' '162.65'
$ws.Range('E28').Value = '  -1.60%  '
Set-TextValue $ws 'D29' '2.517'
$ws.Range('E29').Value = '  -2.01%  '
Set-TextValue $ws 'D30' '132.41'
$ws.Range('E30').Value = '  -2.91%  '
$ws.Range('E31').Value = '  -5.92%  '
$ws.Range('E32').Value = '  -3.43%  '
Set-TextValue $ws 'D33' '1.634'
$ws.Range('E33').Value = '  -1.91%  '
Set-TextValue $ws 'D34' '6.152'
$ws.Range('E34').Value = '  -3.87%  '
Set-TextValue $ws 'D35' '3.963'
$ws.Range('E35').Value = '  -1.36%  '
Set-TextValue $ws 'D36' '6.022'
$ws.Range('E36').Value = '  +1.64%  '
Set-TextValue $ws 'D37' '10.23'
$ws.Range('E37').Value = '  -1.76%  '
Set-TextValue $ws 'D38' '0.02571'
$ws.Range('E38').Value = '  -3.17%  '
Set-TextValue $ws 'D39' '0.06677'
$ws.Range('E39').Value = '  -3.09%  '
$ws.Range('E40').Value = '  -1.99%  '
Set-TextValue $ws 'D41' '0.6837'
$ws.Range('E41').Value = '  -1.30%  '
Set-TextValue $ws 'D42' '0.2226'
$ws.Range('E42').Value = '  -4.23%  '
Set-TextValue $ws 'D43' '1.295'
$ws.Range('E43').Value = '  +1.37%  '
Set-TextValue $ws 'D44' '0.6672'
$ws.Range('E44').Value = '  +3.02%  '
Set-TextValue $ws 'D45' '14.19'
$ws.Range('E45').Value = '  -4.17%  '
Set-TextValue $ws 'D46' '2.292'
$ws.Range('E46').Value = '  -1.86%  '
Set-TextValue $ws 'D47' '3.606'
$ws.Range('E47').Value = '  -3.73%  '
$ws.Range('E48').Value = '  -2.86%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws 'D49' '81.79'
$ws.Range('E49').Value = '  -1.17%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws 'D50' '0.00000000334'
$ws.Range('E50').Value = '  -9.01%  '
Set-TextValue $ws 'D51' '1.160'
$ws.Range('E51').Value = '  -2.68%  '
